$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

function Replace-AndAddSpaceRun($old, $new) {
    # Replace heading text, then append a new, non-bold run containing a single space
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: heading not found -> $old"
        return
    }
    $rng2 = $d.Content
    $rng2.Find.Execute($new, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rng2.Find.Found) {
        Write-Output "WARNING: could not relocate heading -> $new"
        return
    }
    $insStart = $rng2.End
    $rng2.Collapse(0)
    $rng2.InsertAfter(" ")
    $newRng = $d.Range($insStart, $insStart + 1)
    $newRng.Font.Bold = 0
}

Replace-Text "Appendix 7: SWIFT Focus Group Information Sheet and Consent Form: " "ISihlomelo -7: Iphepha loLwazi leQela eliGxilileyo le-SWIFT kunye neFomu yeMvume: "
Replace-Text "Participants" "Abathathi-nxaxheba"
Replace-Text "Uyamenywa ukuba ujoyine ingxoxo yeqela eligxilileyo malunga namava akho ngenkqubo yethu. Oludliwano-ndlebe luyinxalenye yophononongo, olwenziwa ngabaphandi abaphuma kwiDyunivesithi yaseKapa kunye neyoMzantsi Afrika kunye neDyunivesithi yaseOxford eUnited Kingdom. " "You’re invited to join a focus group discussion about your experience with our programme. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. "
Replace-Text "Ngaphambili kokuba ugqibe ekubeni ungathanda na ukuba nodliwano-ndlebe, kubalulekile ukuba wazi kutheni sisenza oluphando nje kwaye ukuthatha inxaxheba kungaquka ntoni. Lonke ulwazi ozakuludinga luchaziwe ngezantsi kodwa ukuba unayo nayiphi na imibuzo malunga nokuthatha inxaxheba okanye ngoluphononongo lwethu, ndicela u-imeyilele iqela lophononongo ku " "Before you decide if you’d like to be interviewed, it’s important for you to know why we’re doing this research and what participating in it would involve. All the information you might need is explained below but if you have any questions about your participation or our study, please email the study team at "
Replace-Text " okanye uthumele umyalezo kuthi ku WhatsApp at +27 XX XXX XXXX. Silapha ukuzokunceda wena!" " or message us on WhatsApp at +27 XX XXX XXXX. We’re here to help you!"
Replace-Text "Yintoni iqela ekugxilwe kulo kwaye eli liza kujongeka njani?" "What is a focus group and what will this one look like?"
Replace-Text "Iqela ekugxilwe kulo yingxoxo yeqela yophando. Njengenxalenye yolu phononongo, uya kubekwa kwiqela labantu abasithandathu ukuya kwabasibhozo. Amanye amalungu eqela onke aya kuba ngabazali kunye nabanonopheli nabo bebesebenza ngenkxaso yobuzali kwi-chatbot. Imodareyitha iya kubuza iqela imibuzo eya kukhokelela kwingxoxo. Kusenokubakho umntu othatha amanqaku/itoliki kwigumbi neqela. Bobabini imodareyitha kunye nomthathi-manqaku/itoliki bayinxalenye yeqela lophononongo. " "A focus group is a group discussion for research. As part of this study, you will be placed in a group of 6 – 8 individuals. The other members of the group will all be fellow parents and caregivers who have also been working through the parenting support on the chatbot. A moderator will ask the group questions that will lead to discussion. There might also be a note-taker/interpreter in the room with the group. Both the moderator and the note-taker/interpreter are part of the study team. "
Replace-Text "Kutheni ndimenyiwe nje kudliwano-ndlebe?" "Why have I been invited to the interview?"
Replace-Text "Umenyiwe kwingxoxo yeqela (kunye nabanye abazali/abanonopheli) kunye nelungu leqela lethu lophando kuba uyinxalenye yophononongo lwethu. Singathanda ukuva malunga namava akho ngenkqubo ye-ParentText. Ukuze ube nodliwano-ndlebe, kufuneka uvume ukuthatha inxaxheba. " "You’ve been invited to a group discussion (along with other parents/caregivers) with a member of our research team because you’re part of our study. We would love to hear about your experience with the ParentText programme. To be interviewed, you need to agree to take part. "
Replace-Text "Hayi, kuxhomekeke kuwe ukuba uyafuna ukujoyina okanye awufuni. Ukuba awukufuni ukwenziwa udliwano-ndlebe, akuzokubakho ziphumo kuwe okanye kusapho lwakho. If you do choose to participate in the groups but don't want to answer some of the questions, you can stop at any time by telling your interviewer or just to stop responding in the group." "Hayi, kuxhomekeke kuwe ukuba uyafuna ukujoyina okanye awufuni. If you don't want to be interviewed, there will be no implications to you or your family. If you do choose to participate in the groups but don't want to answer some of the questions, you can stop at any time by telling your interviewer or just to stop responding in the group."
Replace-Text "Ukuba uthatha isigqibo sokuba ungathanda udliwano-ndlebe, kuya kufuneka ukuba uvumelane ngomlomo kwimibuzo yemvume engezantsi apho umntu okwenza udliwano-ndlebe eya kukubuza khona. Udliwano-ndlebe luya kwenzeka ngobuqu kwaye luya kuba malunga neyure enye ukuya kwiyure emenizuzu enamashumi amahlanu. Ingxoxo iya kuqhutywa kwindawo yabucala elungiswe liqela lophando. " "If you decide you’d like to be interviewed, you’ll need to agree verbally to the consent questions below which the person interviewing you will ask you. Udliwano-ndlebe luya kwenzeka ngobuqu kwaye luya kuba malunga neyure enye ukuya kwiyure emenizuzu enamashumi amahlanu. The discussion will be conducted in a private space arranged by the research team. "
Replace-Text "Ngexesha lodliwano-ndlebe, ilungu leqela lophando liya kukubuza imibuzo malunga neengcinga kunye namava akho okusebenzisa i-chatbot. Sifuna ukubona ukuba abazali bayathanda na ukusebenzisa i-chatbot. Sikwafuna ukwazi ukuba bonwabile na ngemiyalezo kwaye ukuba ukusebenzisa i-chatbot kuyayitshintsha indlela abakhathalela ngayo abantwana babo. Uya kuba nelungelo lokutsiba imibuzo ongafuni ukuyiphendula. Akukho zimpendulo zichanekileyo okanye ezingachanekanga kuba amava akho onke abalulekile kuthi. " "During the interview, a member of the research team will ask you some questions about your thoughts and experiences using the chatbot. We want to see if parents like using the chatbot. We also want to know if they're happy with the messages and if using the chatbot changes how they take care of their kids. You will have the right to skip questions you do not want to answer. There are also no right or wrong answers because your whole experience is important to us. "
Replace-Text "Ukukhusela iinkcukacha zakho (kuquka igama lakho lokwenene, inkcukacha zoqhagamshelwano, kunye naluphi na olunye ulwazi olungakuchaza wena), sizakunika inombolo yokuthatha inxaxheba, kwaye ungazikhethela igama ofuna sikubize ngalo ngexesha lodliwano-ndlebe. Nceda ungabhekisi naliphi na elinye iqela lesithathu ngegama ngexesha lodliwano-ndlebe, ngaphandle kwemvume yabo, ukuze sikwazi ukukhusela iinkcukacha zabo zobuqu. " "To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview. Please also do not refer to any other third parties by name during the interview, without their permission, so that we can protect their personal information too. "
Replace-Text "Sizakushicilela oludliwano-ndlebe ukusinceda sikhumbule ebesixoxe ngako kwaye kamva sikubhale phantsi obekuthethiwe. Uyakunikwa inombolo endaweni yokuba igama lakho lisetyenziswe ukuze naluphi na ulwazi owabelana ngalo kudliwano-ndlebe lweqela lakho lungakwazi ukudityaniswa nawe nguye nabani na ngaphandle kweqela lophando. Sizakuzicima naziphi na inkcukacha zakho esiziqokeleleyo kuwe ekupheleni koluphononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokhelela ekukuchazeni kwixesha lokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kusetyenzwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu agunyazisiweyo kuphela eqela lophando aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa." "We will record the interview to help us remember the discussion and later write down what was said. You will be given a number instead of your name being used so that any information you share in your group interview will not be able to be linked to you by anyone besides the research team. We will delete any personal information we collect from you at the end of the study and, after transcribing your interview, change any data which might lead to identification at the point of transcription. We may use an artificial intelligence (AI) software, Microsoft Transcriber, to transcribe the interviews at first, and then we will check/review these transcriptions. This AI-generated information will be processed and stored securely on password-protected University of Cape Town servers, and in accordance with POPIA. Only authorised members of the research team will be able to access it, and this data will be owned by the Global Parenting Initiative at the University of Cape Town."
Replace-Text "Siyakucela ukuba uhloniphe abanye abantu kwiqela, kwaye ungaxoxi ngento ethethwa ngabanye, ngaphandle kwengxoxo yeqela. Siza kuqinisekisa ukuba iingxelo zethu zibhaliwe ukuze kungabikho mntu unokukuchaza kule ngxelo. Nceda ukhumbule, nangona kunjalo, ukuba sinokuqinisekisa oku kuphela kwiqela lophando." "We ask you to respect the other people in the group, and not to discuss what is said by others, outside of the group discussion. We will make sure that our reports are written so that no-one can identify you from the report. Please remember, though, that we can only guarantee this for the research team."
Replace-Text "Ingaba ikhona into endiyifumanayo ngokwenziwa oludliwano-ndlebe? " "Do I get anything for being interviewed? "
Replace-Text "Njengombulelo ngokuthatha inxaxheba kwingxoxo, siza kukunika ivawutsha ye-R120 yakwaShoprite emva koko." "As a thank you for taking part in the discussion, we'll give you a R120 Shoprite voucher afterwards. "
Replace-Text "Sizakuqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcine ngokukhuselekileyo. Ulwazi lwakho, olufana nefomu yakho yemvume kunye noshicilelo lodliwano-ndlebe, kunye nalo naluphi na ulwazi olunikeza nge-imeyile okanye nge-WhatsApp, luya kugcinwa likhuselekile kwiiseva ezikhuselekileyo kwiDyunivesithi yaseKapa. " "We only collect what’s needed for the study and store it securely. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. "
Replace-Text "Ushicilelo lodliwano-ndlebe luzakucinywa emva kokuba sibhale phantsi amanqaku ethu. Naziphi na iinkcukacha ezichaza wena zizakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kucinywa xa uphononongo liphelile. " "Interview recordings will be deleted after we have written our notes. Any details that identify you will be kept separate and only authorised staff can access them. All data will be kept for five years after the study, but personal information will be deleted when the study ends. "
Replace-Text "Iikomiti zokuziphatha kunye nabahloli banokujonga ulwazi. Iinkcukacha zakho ziyakuhlala ziyimfihlo ngaphandle kokuba umthetho uthetha enye into. Emva kophononongo, singabelana ngolwazi nabanye abaphandi kodwa ngaphandle kweenkcukacha zakho. Unelungelo lokubona, ulungise, okanye ucele ukuba kucinywe ulwazi lwakho." "Ethics committees and monitors may check the information. Your information will stay private unless the law says otherwise. After the study, we may share the information with other researchers but without your details. You have the right to see, correct, or ask us to delete your personal information."
Replace-Text "Unelungelo lokucela ukufikelela kwidatha yakho, ukulungisa naziphi na iimpazamo kwidatha yakho, kwaye usicele ukuba siyicime okanye siyidlulisele kwenye indawo. Nceda u-imeyilele iqela lophononongo phambi kwe [*umhla oza kumiselwa] ukuba ufuna ukwenza nayiphi na kwezi." "You have the right to request access to your data, to correct any mistakes in your data, and to request us to delete it or transfer it somewhere else. Please email the study team before [*date to be determined] if you would like to do any of these."
Replace-Text "Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. Siceba ukwabelana ngeziphumo kwiingxelo nakwii-nkomfa ukuze nabanye bafunde kolu phononongo." "Your participation and what you tell us will help us understand how to support families like yours. Siceba ukwabelana ngeziphumo kwiingxelo nakwii-nkomfa ukuze nabanye bafunde kolu phononongo."
Replace-Text "Abaphononongi abaziintloko kolu phononongo nguNjinga Cathy Ward no Cindee Bruyns ze Co-investigator ngu Carly Katzef bonke basuka kwiDyunivesithi yaseKapa." "The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town."
Replace-Text "Ingaba ikhona imingcipheko ekwenziweni udliwano-ndlebe?   " "Are there any risks in being interviewed?   "
Replace-Text "Asilindelanga nayiphi na imingcipheko kuwe ukuba unodliwano-ndlebe. Ukuba nayiphi na imibuzo ikwenza ungakhululeki, akunyanzelekanga uyiphendule. Ukuba uyacaphuka ngexesha lodliwano-ndlebe, unokwazisa umenzi wodliwano-ndlebe lwakho. Khumbula, ungayeka ukuthatha inxaxheba nanini na ngaphandle kokunikeza isizathu. Siyakhathala ngempilo-ntle yakho." "We don’t expect any risks to you if you are interviewed. If any questions make you uncomfortable, you don’t have to answer them. If you become upset during the interview, you can let your interviewer know. Remember, you can stop participating anytime without giving a reason. We care about your well-being."
Replace-Text "Sifuna nokuqinisekisa ukuba ukhuselekile. Ukuba siyaqaphelisisa ukuba wena okanye usapho lwakho lukweyona inkulu ingozi, singanithumela ukuze nifumane inkxaso okanye singadinga ukucela uncedo kwezinye indawo ezingaphandle koluphononongo, njengoo nontlalontle okanye uncedo lwezempilo." "We also want to make sure you're safe. If we notice that you or your family are in serious danger, we might refer you for support or could need to ask for help from other places outside of this study, like social or medical services."
Replace-Text "Ngubani obhatalela oluphononongo?" "Who pays for the study?"
Replace-Text "Olu phononongo luyinxalenye ye Global Parenting Initiative, luxhaswe ngokwezimali ngu LEGO Foundation, Oak Fundation, i-World Childhood Foundation, i-Human Safety Net kunye ne UK Research kunye ne Innovaion Global Challenges Research Fund. " "This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. "
Replace-Text "Ukhuseleko lwedatha" "Data protection"
Replace-Text "Idyunivesithi yaseKapa iqinisekisa ukuba iinkcukacha zakho zobuqu zisetyenziswa ngokukhuselekileyo nangokuchanekileyo, nje kuphando kuphela. Uphononongo lulandela imithetho yokukhuselwa kwedatha efana ne-GDPR (General Data Protection Regulation) e-UK kunye ne-POPIA (uMthetho woKhuselo loLwazi loMntu) eMzantsi Afrika. Nayiphi na idatha ethi ithunyelwe ngaphesheya kwemida izakuthobelana ne POPIA. " "The University Cape Town makes sure your personal information is used safely and correctly, just for research. The study follows data protection laws like GDPR (General Data Protection Regulation) in the UK and POPIA (Protection of Personal Information Act) in South Africa. Any data that is transferred across borders will comply with POPIA. "
Replace-Text "Who has approved this study?" "Ngubani ogunyazise oluphononongo?"

Replace-AndAddSpaceRun "Kwenzeka ntoni ngeenkcukacha zam ukuba ndiyavuma ukuba noludliwano-ndlebe?" "What happens to my information if I agree to be interviewed?"
Replace-AndAddSpaceRun "Ngobani amanye amalungu eqela lophononongo?" "Who are some of the study team members?"

Write-Output "All replacements complete"
